$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (9 rows) replacing previous 4 rows.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.
$data = @(
    @("ECs","Selplg","Sele","ECs",3,1,1.045864666666667,3.137594,0.01599546767026659,0.01599546767026659,3,1,4.699506666666667,14.09852,0.9660495246229048,0.9660495246229047,4.915047973431111,44.23543176088,0.01545241393898208,0.01545241393898208),
    @("ECs","Selplg","Sele","FAPs",3,1,1.045864666666667,3.137594,0.01599546767026659,0.01599546767026659,1,0.3333333333333333,0.1651576666666667,0.495473,0.03395047537709522,0.03395047537709522,0.1727325679957778,1.554593111962,0.0005430537312845084,0.0005430537312845084),
    @("FAPs","Selplg","Sele","ECs",3,1,2.006702333333334,6.020107,0.03069053130839924,0.03069053130839923,3,1,4.699506666666667,14.09852,0.9660495246229048,0.9660495246229047,9.430510993515558,84.87459894164,0.02964857318090346,0.02964857318090345),
    @("FAPs","Selplg","Sele","FAPs",3,1,2.006702333333334,6.020107,0.03069053130839924,0.03069053130839923,1,0.3333333333333333,0.1651576666666667,0.495473,0.03395047537709522,0.03395047537709522,0.331422275067889,2.982800475611,0.001041958127495778,0.001041958127495778),
    @("M2","Selplg","Sele","ECs",3,1,59.82593033333333,179.477791,0.914978548362651,0.914978548362651,3,1,4.699506666666667,14.09852,0.9660495246229048,0.9660495246229047,281.1523584410356,2530.37122596932,0.8839145916858946,0.8839145916858945),
    @("M2","Selplg","Sele","FAPs",3,1,59.82593033333333,179.477791,0.914978548362651,0.914978548362651,1,0.3333333333333333,0.1651576666666667,0.495473,0.03395047537709522,0.03395047537709522,9.880711060015889,88.926399540143,0.03106395667675651,0.03106395667675651),
    @("sCs","Selplg","Sele","ECs",3,1,2.506566,7.519698,0.03833545265868316,0.03833545265868316,3,1,4.699506666666667,14.09852,0.9660495246229048,0.9660495246229047,11.77962362744,106.01661264696,0.03703394581712473,0.03703394581712473),
    @("sCs","Selplg","Sele","FAPs",3,1,2.506566,7.519698,0.03833545265868316,0.03833545265868316,1,0.3333333333333333,0.1651576666666667,0.495473,0.03395047537709522,0.03395047537709522,0.413978591906,3.725807327154,0.001301506841558422,0.001301506841558422)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
